# Applies the "añadí calculo de n de cada subsubmuestra" edit:
#  - Sheet1: mark J13 with the 17-decimal number format (new custom numFmt
#    + cellXf), nudge the H/I column widths, and move the selection/zoom.
#  - Add a new "Sheet2" (placed right after Sheet1) holding the per-
#    subsample n/estudio tables, and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New Sheet2, inserted immediately after Sheet1 --------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Block 1: "5p"
$ws2.Range("A1").Value = "5p"
$ws2.Range("A2").Value = "estudio"
$ws2.Range("B2").Value = "n"

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = 145
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = 116
$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = 202
$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = 281
$ws2.Range("A7").Value = 6
$ws2.Range("B7").Value = 235
$ws2.Range("A8").Value = 17
$ws2.Range("B8").Value = 280
$ws2.Range("A9").Value = 18
$ws2.Range("B9").Value = 293

$ws2.Range("A10").Value = 20
$ws2.Range("B10").Value = 233
$ws2.Range("T10").Value = 235
$ws2.Range("X10").Value = 280
$ws2.Range("AB10").Value = 293
$ws2.Range("AF10").Value = 233

# Block 2: "7p"
$ws2.Range("A11").Value = "7p"
$ws2.Range("A12").Value = "estudio"
$ws2.Range("B12").Value = "n"

$ws2.Range("A13").Value = 11
$ws2.Range("B13").Value = 254
$ws2.Range("T13").Value = 251
$ws2.Range("X13").Value = 254
$ws2.Range("AB13").Value = 156

$ws2.Range("A14").Value = 12
$ws2.Range("B14").Value = 175
$ws2.Range("A15").Value = 13
$ws2.Range("B15").Value = 294
$ws2.Range("A16").Value = 14
$ws2.Range("B16").Value = 293
$ws2.Range("A17").Value = 15
$ws2.Range("B17").Value = 251
$ws2.Range("A18").Value = 16
$ws2.Range("B18").Value = 254
$ws2.Range("A19").Value = 19
$ws2.Range("B19").Value = 156

# --- Sheet1 tweaks ------------------------------------------------------
# New "n" column (J) gets a 17-decimal custom number format; the cell
# itself stays empty (style only), same as the source edit.
$ws1.Range("J13").NumberFormat = "0.00000000000000000"

# Column width nudges following the new J column (H/I shrink-to-fit changed
# slightly, I lost its auto bestFit, new J column got an explicit width).
$ws1.Columns.Item(8).ColumnWidth = 5.42
$ws1.Columns.Item(9).ColumnWidth = 7.09
$ws1.Columns.Item(10).ColumnWidth = 25.76

# Selection / zoom moved on Sheet1 before the user switched to Sheet2.
$ws1.Range("E6").Select()
$excel.ActiveWindow.Zoom = 190

# Sheet2 ends up the active/selected tab.
$ws2.Select()
